$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 142, shifting existing rows 142:198 down to 143:199
$ws.Rows.Item(142).Insert()

# Populate the newly inserted row 142 with the new record's data.
# Columns A,B,C,E,F,G,H,I,R are identical for every row in this block,
# so copy them from the row above (row 141) which is untouched by the insert.
$ws.Range("A142").Value = $ws.Range("A141").Value2
$ws.Range("B142").Value = $ws.Range("B141").Value2
$ws.Range("C142").Value = $ws.Range("C141").Value2
$ws.Range("E142").Value = $ws.Range("E141").Value2
$ws.Range("F142").Value = $ws.Range("F141").Value2
$ws.Range("G142").Value = $ws.Range("G141").Value2
$ws.Range("H142").Value = $ws.Range("H141").Value2
$ws.Range("I142").Value = $ws.Range("I141").Value2
$ws.Range("R142").Value = $ws.Range("R141").Value2

# New record-specific values
$ws.Range("D142").Value = 45141
$ws.Range("D142").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("J142").Value = 180
$ws.Range("K142").Value = 24000
$ws.Range("L142").Value = 25000
$ws.Range("M142").Value = 24556
$ws.Range("N142").Value = "$/malla 15 kilos"
$ws.Range("O142").Value = "Provincia de Quillota"
$ws.Range("P142").Value = 1637
$ws.Range("Q142").Value = 15
